$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.5490618995054
$ws.Range("D2").Value = 3.489273697077929
$ws.Range("E2").Value = 25.80812753231834
$ws.Range("F2").Value = 24.39104713117687
$ws.Range("G2").Value = 30.93928248280366
$ws.Range("H2").Value = 13.48434404676267
$ws.Range("I2").Value = 24.64241205178492
$ws.Range("L2").Value = 10.24440483804235
$ws.Range("M2").Value = 15.80007476556411
$ws.Range("B3").Value = 16.07345499580249
$ws.Range("D3").Value = 3.527311135262578
$ws.Range("E3").Value = 24.98561795620806
$ws.Range("F3").Value = 24.0189881537857
$ws.Range("G3").Value = 30.08631479192796
$ws.Range("H3").Value = 13.46243202893969
$ws.Range("I3").Value = 24.79833153369635
$ws.Range("L3").Value = 10.05779004209083
$ws.Range("M3").Value = 15.55856122571748
$ws.Range("B4").Value = 15.77461083674064
$ws.Range("D4").Value = 3.551683636425415
$ws.Range("E4").Value = 24.46528320796558
$ws.Range("F4").Value = 23.80057920529715
$ws.Range("G4").Value = 29.56972175102106
$ws.Range("H4").Value = 13.45439511741066
$ws.Range("I4").Value = 24.90250789127501
$ws.Range("L4").Value = 9.94255819574655
$ws.Range("M4").Value = 15.40950760476522
$ws.Range("B5").Value = 15.65127198566043
$ws.Range("D5").Value = 3.561872532496356
$ws.Range("E5").Value = 24.2496457394453
$ws.Range("F5").Value = 23.71422809194632
$ws.Range("G5").Value = 29.36142230557179
$ws.Range("H5").Value = 13.45247851268861
$ws.Range("I5").Value = 24.94707638217177
$ws.Range("L5").Value = 9.895491934239759
$ws.Range("M5").Value = 15.34863674177978
$ws.Range("B6").Value = 15.63070242244895
$ws.Range("D6").Value = 3.563579938964909
$ws.Range("E6").Value = 24.21362991331909
$ws.Range("F6").Value = 23.7000533431437
$ws.Range("G6").Value = 29.32698068811776
$ws.Range("H6").Value = 13.45224217003122
$ws.Range("I6").Value = 24.95460453701036
$ws.Range("L6").Value = 9.887671703304072
$ws.Range("M6").Value = 15.33852312892611
$ws.Range("B7").Value = 15.77295353551562
$ws.Range("D7").Value = 3.551820005859756
$ws.Range("E7").Value = 24.46238925746249
$ws.Range("F7").Value = 23.79940374129107
$ws.Range("G7").Value = 29.56690301415482
$ws.Range("H7").Value = 13.45436377519813
$ws.Range("I7").Value = 24.90310039983692
$ws.Range("L7").Value = 9.941923808914099
$ws.Range("M7").Value = 15.40868712463044
$ws.Range("B8").Value = 16.38657892323738
$ws.Range("D8").Value = 3.50217850602652
$ws.Range("E8").Value = 25.52784866471774
$ws.Range("F8").Value = 24.26075752052503
$ws.Range("G8").Value = 30.64395006388795
$ws.Range("H8").Value = 13.47566276627949
$ws.Range("I8").Value = 24.69441698984819
$ws.Range("L8").Value = 10.18022749010127
$ws.Range("M8").Value = 15.71699601995794
$ws.Range("B9").Value = 17.52946041729144
$ws.Range("D9").Value = 3.41285404322861
$ws.Range("E9").Value = 27.48577194905711
$ws.Range("F9").Value = 25.23910722232923
$ws.Range("G9").Value = 32.79457078547831
$ws.Range("H9").Value = 13.56049428509368
$ws.Range("I9").Value = 24.35247772302744
$ws.Range("L9").Value = 10.64004050805717
$ws.Range("M9").Value = 16.31287931479059
$ws.Range("B10").Value = 18.3246512515449
$ws.Range("D10").Value = 3.352047235009518
$ws.Range("E10").Value = 28.83236289625372
$ws.Range("F10").Value = 25.99440311758435
$ws.Range("G10").Value = 34.37508552015308
$ws.Range("H10").Value = 13.64906456288469
$ws.Range("I10").Value = 24.14271104522942
$ws.Range("L10").Value = 10.97027578728088
$ws.Range("M10").Value = 16.74191447103779
$ws.Range("B11").Value = 18.6753611271996
$ws.Range("D11").Value = 3.32541605064
$ws.Range("E11").Value = 29.42294076069175
$ws.Range("F11").Value = 26.34410190654116
$ws.Range("G11").Value = 35.08975147781585
$ws.Range("H11").Value = 13.69501606428515
$ws.Range("I11").Value = 24.05638062274643
$ws.Range("L11").Value = 11.11826309488901
$ws.Range("M11").Value = 16.93449591609414
$ws.Range("B12").Value = 18.80647971455932
$ws.Range("D12").Value = 3.315478496162535
$ws.Range("E12").Value = 29.64326786212676
$ws.Range("F12").Value = 26.47725003870147
$ws.Range("G12").Value = 35.35941983223381
$ws.Range("H12").Value = 13.71322437592541
$ws.Range("I12").Value = 24.02500533940144
$ws.Range("L12").Value = 11.17393352305642
$ws.Range("M12").Value = 17.00699396449058
$ws.Range("B13").Value = 18.77831755195422
$ws.Range("D13").Value = 3.317612200566635
$ws.Range("E13").Value = 29.59596596795404
$ws.Range("F13").Value = 26.44854439583167
$ws.Range("G13").Value = 35.30138974301872
$ws.Range("H13").Value = 13.70926709970632
$ws.Range("I13").Value = 24.03170389389456
$ws.Range("L13").Value = 11.16196104331519
$ws.Range("M13").Value = 16.99140013122503
$ws.Range("B14").Value = 18.6861827021125
$ws.Range("D14").Value = 3.324595539924815
$ws.Range("E14").Value = 29.44113437238162
$ws.Range("F14").Value = 26.35504243083873
$ws.Range("G14").Value = 35.11195834508018
$ws.Range("H14").Value = 13.6964979335239
$ws.Range("I14").Value = 24.05377293198295
$ws.Range("L14").Value = 11.12285079438294
$ws.Range("M14").Value = 16.94046927181568
$ws.Range("B15").Value = 18.62952481793388
$ws.Range("D15").Value = 3.328892162646301
$ws.Range("E15").Value = 29.3458601097958
$ws.Range("F15").Value = 26.29785963780665
$ws.Range("G15").Value = 34.99579149972707
$ws.Range("H15").Value = 13.68878138469048
$ws.Range("I15").Value = 24.06746251474727
$ws.Range("L15").Value = 11.09884520451052
$ws.Range("M15").Value = 16.90921524103866
$ws.Range("B16").Value = 18.30150069155097
$ws.Range("D16").Value = 3.353808285913098
$ws.Range("E16").Value = 28.79331151276849
$ws.Range("F16").Value = 25.97165923678852
$ws.Range("G16").Value = 34.32826391917121
$ws.Range("H16").Value = 13.64617491282121
$ws.Range("I16").Value = 24.14853665369326
$ws.Range("L16").Value = 10.96055532594449
$ws.Range("M16").Value = 16.72927187629846
$ws.Range("B17").Value = 18.09736839567796
$ws.Range("D17").Value = 3.369356625128556
$ws.Range("E17").Value = 28.44859591254543
$ws.Range("F17").Value = 25.7730008356388
$ws.Range("G17").Value = 33.91740156507191
$ws.Range("H17").Value = 13.62148294900381
$ws.Range("I17").Value = 24.20060812067803
$ws.Range("L17").Value = 10.87511120191157
$ws.Range("M17").Value = 16.61817871488985
$ws.Range("B18").Value = 17.97892444648144
$ws.Range("D18").Value = 3.378396641686686
$ws.Range("E18").Value = 28.24826226349725
$ws.Range("F18").Value = 25.65932209237955
$ws.Range("G18").Value = 33.68069657593067
$ws.Range("H18").Value = 13.60781407578529
$ws.Range("I18").Value = 24.23141365696677
$ws.Range("L18").Value = 10.82575798115555
$ws.Range("M18").Value = 16.55404040317215
$ws.Range("B19").Value = 17.93864731925655
$ws.Range("D19").Value = 3.381474131374723
$ws.Range("E19").Value = 28.18008336524239
$ws.Range("F19").Value = 25.62093750452975
$ws.Range("G19").Value = 33.60049649813121
$ws.Range("H19").Value = 13.60327779173991
$ws.Range("I19").Value = 24.24199054014851
$ws.Range("L19").Value = 10.8090135839741
$ws.Range("M19").Value = 16.53228472208465
$ws.Range("B20").Value = 18.11920625462853
$ws.Range("D20").Value = 3.367691443876268
$ws.Range("E20").Value = 28.48550606061758
$ws.Range("F20").Value = 25.79408900952085
$ws.Range("G20").Value = 33.96118124555768
$ws.Range("H20").Value = 13.62405629842692
$ws.Range("I20").Value = 24.19497643132788
$ws.Range("L20").Value = 10.8842287548692
$ws.Range("M20").Value = 16.63003007276858
$ws.Range("B21").Value = 18.71329150389308
$ws.Range("D21").Value = 3.322540379929638
$ws.Range("E21").Value = 29.48670315238403
$ws.Range("F21").Value = 26.38248776877538
$ws.Range("G21").Value = 35.16762756572219
$ws.Range("H21").Value = 13.70022669302587
$ws.Range("I21").Value = 24.04725493143858
$ws.Range("L21").Value = 11.13434880502971
$ws.Range("M21").Value = 16.9554409603718
$ws.Range("B22").Value = 19.09168370331417
$ws.Range("D22").Value = 3.293888466833524
$ws.Range("E22").Value = 30.12168182363142
$ws.Range("F22").Value = 26.7711955933218
$ws.Range("G22").Value = 35.95037538914686
$ws.Range("H22").Value = 13.75471078235804
$ws.Range("I22").Value = 23.95838608354059
$ws.Range("L22").Value = 11.29564494633693
$ws.Range("M22").Value = 17.16559469258155
$ws.Range("B23").Value = 18.89066342383715
$ws.Range("D23").Value = 3.309102464895794
$ws.Range("E23").Value = 29.7845981757896
$ws.Range("F23").Value = 26.56340487941651
$ws.Range("G23").Value = 35.53323792664368
$ws.Range("H23").Value = 13.72520393435462
$ws.Range("I23").Value = 24.00511188963019
$ws.Range("L23").Value = 11.20977181454692
$ws.Range("M23").Value = 17.05368026076258
$ws.Range("B24").Value = 18.10933673206002
$ws.Range("D24").Value = 3.36844395743353
$ws.Range("E24").Value = 28.46882567433809
$ws.Range("F24").Value = 25.78455337541636
$ws.Range("G24").Value = 33.94138996700541
$ws.Range("H24").Value = 13.62289124561298
$ws.Range("I24").Value = 24.19751981212374
$ws.Range("L24").Value = 10.88010742275716
$ws.Range("M24").Value = 16.62467290862767
$ws.Range("B25").Value = 17.22758512756725
$ws.Range("D25").Value = 3.436167158013383
$ws.Range("E25").Value = 26.97147146860318
$ws.Range("F25").Value = 24.96740797781834
$ws.Range("G25").Value = 32.21110098141605
$ws.Range("H25").Value = 13.53292841217218
$ws.Range("I25").Value = 24.43773319940622
$ws.Range("L25").Value = 10.51677111453676
$ws.Range("M25").Value = 16.15297237447164
